# "auf den neuesten Stand gebracht" - bring the project inventory up to date:
# refresh start/end dates for several projects and remove a batch of
# projects (rows 12-21) that no longer belong on the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

# --- refreshed Start (C) / Ende (D) dates for still-active projects ---

# Row 3
$ws.Range("D3").Value = 43110

# Row 5
$ws.Range("C5").Value = 41821
$ws.Range("D5").Value = 42551

# Row 6
$ws.Range("C6").Value = 41760
$ws.Range("D6").Value = 42004

# Row 7
$ws.Range("C7").Value = 41791
$ws.Range("D7").Value = 43220

# Row 8
$ws.Range("C8").Value = 41821
$ws.Range("D8").Value = 42886

# Row 9 - still running, no end date anymore
$ws.Range("C9").Value = 42217
$ws.Range("D9").ClearContents()

# Row 11
$ws.Range("C11").Value = 41927
$ws.Range("D11").Value = 42613

# --- remove the now-obsolete projects in rows 12-21 ---
$ws.Range("A12:G21").ClearContents()

# --- move the active selection to where work continues ---
$ws.Range("A18").Select()
